$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.438.41"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.905.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "483.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("E7").Value = "  -2.51%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.44%  "

$ws.Range("E10").Value = "  +7.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000351"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.37%  "

$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.524.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.916.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.75%  "

$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.21%  "

$ws.Range("E19").Value = "  -3.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.502.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("E22").Value = "  -2.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "712.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.22%  "

$ws.Range("E31").Value = "  -4.32%  "

$ws.Range("E32").Value = "  -4.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0886"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +31.24%  "

$ws.Range("E35").Value = "  -4.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.92%  "

$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.339"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.08%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("E47").Value = "  -1.62%  "

$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "148.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.02%  "

$ws.Range("E51").Value = "  -3.52%  "
